# Nexial "#system" sheet update:
#  - new category "aws.vision" -> new column E (pushes former E..AD to F..AE)
#  - new category name inserted alphabetically into column A ("target" list)
#  - new "web" function screenshot(file,locator) inserted alphabetically
#  - 3 renamed function labels (NotContains -> NotContain)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert a brand-new column before column E. Everything that used to
#    live in E..AD (base, csv, desktop, ... xml) shifts right to F..AE.
# ---------------------------------------------------------------------
$ws.Range("E1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2) Populate the new column E with the "aws.vision" category.
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "aws.vision"
$ws.Range("E2").Value = "saveText(profile,image,var)"

# ---------------------------------------------------------------------
# 3) Insert "aws.vision" into the alphabetical category index in column A
#    (A2:A30 -> A2:A31). Only column A shifts -- shift existing rows 30
#    down to 5 manually so every other column is untouched.
# ---------------------------------------------------------------------
for ($r = 30; $r -ge 5; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Text
}
$ws.Range("A5").Value = "aws.vision"

# ---------------------------------------------------------------------
# 4) Insert the new "screenshot(file,locator)" command into the "web"
#    list (now column Z after the column insert above), alphabetically
#    between saveValues(...) and scrollElement(...) -> row 107, pushing
#    rows 107-134 down to 108-135.
# ---------------------------------------------------------------------
for ($r = 134; $r -ge 107; $r--) {
    $ws.Cells.Item($r + 1, 26).Value = $ws.Cells.Item($r, 26).Text
}
$ws.Range("Z107").Value = "screenshot(file,locator)"

# ---------------------------------------------------------------------
# 5) Rename three grammatically-displeasing labels (still same cells,
#    just shifted one column right by step 1's insert).
# ---------------------------------------------------------------------
$ws.Range("F11").Value = "assertNotContain(text,substring)"
$ws.Range("Z5").Value = "assertAttributeNotContain(locator,attrName,contains)"
$ws.Range("Z41").Value = "assertTextNotContain(locator,text)"

# ---------------------------------------------------------------------
# 6) Fix up every defined name whose column shifted right by one letter.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("csv").RefersTo = "='#system'!`$G`$2:`$G`$6"
$wb.Names.Item("desktop").RefersTo = "='#system'!`$H`$2:`$H`$98"
$wb.Names.Item("excel").RefersTo = "='#system'!`$I`$2:`$I`$14"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$5"
$wb.Names.Item("image").RefersTo = "='#system'!`$K`$2:`$K`$7"
$wb.Names.Item("io").RefersTo = "='#system'!`$L`$2:`$L`$29"
$wb.Names.Item("jms").RefersTo = "='#system'!`$M`$2:`$M`$4"
$wb.Names.Item("json").RefersTo = "='#system'!`$N`$2:`$N`$18"
$wb.Names.Item("macro").RefersTo = "='#system'!`$P`$2:`$P`$4"
$wb.Names.Item("mail").RefersTo = "='#system'!`$Q`$2:`$Q`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$S`$2:`$S`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$T`$2:`$T`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$U`$2:`$U`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$V`$2:`$V`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$W`$2:`$W`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$X`$2:`$X`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$Y`$2:`$Y`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$135"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Item("localdb").RefersTo = "='#system'!`$O`$2:`$O`$12"

# ---------------------------------------------------------------------
# 7) Register the brand-new "aws.vision" defined name.
# ---------------------------------------------------------------------
$wb.Names.Add("aws.vision", "='#system'!`$E`$2:`$E`$2")
